# Add a new "18-jun" column (E) to the "Prix Spot" sheet, mirroring the
# existing day columns (B = 15-jun, C = 16-jun, D = 17-jun).
#
# Row 1 gets the new date header with the same (bold / centered / bordered)
# style as the other header cells; rows 2-25 get a placeholder "-" value,
# exactly like how a not-yet-published day is represented in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Header cell: new date label, formatted like the other headers (B1:D1).
$ws.Range("E1").Value = "18-jun"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows: placeholder "-" (no price published yet), unformatted like
# the existing numeric cells in B2:D25.
for ($row = 2; $row -le 25; $row++) {
    $ws.Cells.Item($row, 5).Value = "-"
}
